$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$origStyle = $ws.Range('B2').Style

$ws.Range('D2').Value = '26.794.97'
$ws.Range('E2').Value = '  -2.43%  '
$ws.Range('D3').Value = '1.569.57'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.74'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.00'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0585'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '1.790.81'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '1.563.41'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.74'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('D16').Value = '26.798.22'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.45'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.42'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.53'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').Value = '0.0₃0678'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.33'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.73'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.97'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0466'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.11'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.17'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').Value = '1.396.75'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.94'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('E37').Value = '  -2.16%  '
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.529'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -2.78%  '
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.988'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.32'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.18'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.43'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '1.704.12'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.96'
$ws.Range('D48').Style = $origStyle
$ws.Range('D49').Value = '0.0₇0984'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0955'
$ws.Range('D50').Style = $origStyle
$ws.Range('E51').Value = '  -0.78%  '